$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all Fitness values (column C, rows 2-12) to 4868
$ws.Range("C2:C12").Value = 4868
